{"js": "// Replace each two-digit multiplication expression with its new value.\n// The mapping below reproduces the OOXML diff exactly, in document order:\n// every \"AA\u00d7BB=CCCC\" text run in the table is swapped for a freshly\n// generated (but equally formatted) multiplication fact.\nconst replacements = [\n  [\"65\u00d723=1495\", \"84\u00d767=5628\"],\n  [\"88\u00d726=2288\", \"39\u00d797=3783\"],\n  [\"69\u00d729=2001\", \"17\u00d795=1615\"],\n  [\"63\u00d795=5985\", \"94\u00d729=2726\"],\n  [\"81\u00d753=4293\", \"42\u00d771=2982\"],\n  [\"19\u00d753=1007\", \"88\u00d721=1848\"],\n  [\"80\u00d775=6000\", \"47\u00d716=752\"],\n  [\"36\u00d761=2196\", \"31\u00d777=2387\"],\n  [\"63\u00d767=4221\", \"14\u00d787=1218\"],\n  [\"35\u00d795=3325\", \"65\u00d711=715\"],\n  [\"51\u00d719=969\", \"90\u00d791=8190\"],\n  [\"58\u00d716=928\", \"68\u00d758=3944\"],\n  [\"55\u00d753=2915\", \"83\u00d792=7636\"],\n  [\"11\u00d789=979\", \"37\u00d765=2405\"],\n  [\"39\u00d777=3003\", \"57\u00d739=2223\"],\n  [\"73\u00d737=2701\", \"11\u00d753=583\"],\n  [\"16\u00d756=896\", \"54\u00d780=4320\"],\n  [\"44\u00d797=4268\", \"22\u00d714=308\"],\n  [\"66\u00d772=4752\", \"11\u00d736=396\"],\n  [\"91\u00d713=1183\", \"56\u00d738=2128\"],\n  [\"13\u00d739=507\", \"62\u00d752=3224\"],\n  [\"84\u00d775=6300\", \"29\u00d786=2494\"],\n  [\"40\u00d718=720\", \"90\u00d755=4950\"],\n  [\"65\u00d797=6305\", \"46\u00d719=874\"],\n  [\"26\u00d723=598\", \"92\u00d750=4600\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  // Full-text search: every expression in the sheet is unique, so this\n  // locates exactly the one run that needs to change.\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    // Replacing the matched range's text keeps the existing run\n    // formatting (rFonts/sz) intact, matching the original diff.\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression with its new value.\n# The mapping below reproduces the OOXML diff exactly, in document order:\n# every \"AA\u00d7BB=CCCC\" text run in the table is swapped for a freshly\n# generated (but equally formatted) multiplication fact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"65\u00d723=1495\", \"84\u00d767=5628\"),\n    @(\"88\u00d726=2288\", \"39\u00d797=3783\"),\n    @(\"69\u00d729=2001\", \"17\u00d795=1615\"),\n    @(\"63\u00d795=5985\", \"94\u00d729=2726\"),\n    @(\"81\u00d753=4293\", \"42\u00d771=2982\"),\n    @(\"19\u00d753=1007\", \"88\u00d721=1848\"),\n    @(\"80\u00d775=6000\", \"47\u00d716=752\"),\n    @(\"36\u00d761=2196\", \"31\u00d777=2387\"),\n    @(\"63\u00d767=4221\", \"14\u00d787=1218\"),\n    @(\"35\u00d795=3325\", \"65\u00d711=715\"),\n    @(\"51\u00d719=969\", \"90\u00d791=8190\"),\n    @(\"58\u00d716=928\", \"68\u00d758=3944\"),\n    @(\"55\u00d753=2915\", \"83\u00d792=7636\"),\n    @(\"11\u00d789=979\", \"37\u00d765=2405\"),\n    @(\"39\u00d777=3003\", \"57\u00d739=2223\"),\n    @(\"73\u00d737=2701\", \"11\u00d753=583\"),\n    @(\"16\u00d756=896\", \"54\u00d780=4320\"),\n    @(\"44\u00d797=4268\", \"22\u00d714=308\"),\n    @(\"66\u00d772=4752\", \"11\u00d736=396\"),\n    @(\"91\u00d713=1183\", \"56\u00d738=2128\"),\n    @(\"13\u00d739=507\", \"62\u00d752=3224\"),\n    @(\"84\u00d775=6300\", \"29\u00d786=2494\"),\n    @(\"40\u00d718=720\", \"90\u00d755=4950\"),\n    @(\"65\u00d797=6305\", \"46\u00d719=874\"),\n    @(\"26\u00d723=598\", \"92\u00d750=4600\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
